# Updated cryptos list on Wed Jun 14 21:27:22 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.126.05"
$ws.Range("E2").Value = "'  -2.77%  "

# Row 3
$ws.Range("D3").Value = "'1.649.12"
$ws.Range("E3").Value = "'  -5.06%  "

# Row 4
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "'  +0.34%  "

# Row 5
$ws.Range("D5").Value = "'235.54"
$ws.Range("E5").Value = "'  -2.15%  "

# Row 6
$ws.Range("E6").Value = "'  +0.22%  "

# Row 7
$ws.Range("D7").Value = "'0.4780"
$ws.Range("E7").Value = "'  -8.43%  "

# Row 8
$ws.Range("D8").Value = "'0.2608"
$ws.Range("E8").Value = "'  -4.74%  "

# Row 9
$ws.Range("D9").Value = "'0.05957"
$ws.Range("E9").Value = "'  -3.32%  "

# Row 10
$ws.Range("D10").Value = "'0.07070"
$ws.Range("E10").Value = "'  -1.59%  "

# Row 11
$ws.Range("D11").Value = "'1.665.45"
$ws.Range("E11").Value = "'  -4.42%  "

# Row 12
$ws.Range("B12").Value = "'Solana"
$ws.Range("C12").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'14.33"
$ws.Range("E12").Value = "'  -3.94%  "

# Row 13
$ws.Range("B13").Value = "'Polygon"
$ws.Range("C13").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.6165"
$ws.Range("E13").Value = "'  -3.85%  "

# Row 14
$ws.Range("D14").Value = "'4.563"
$ws.Range("E14").Value = "'  -1.22%  "

# Row 15
$ws.Range("D15").Value = "'72.91"
$ws.Range("E15").Value = "'  -5.81%  "

# Row 16
$ws.Range("E16").Value = "'  +0.25%  "

# Row 17
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "'  +0.33%  "

# Row 18
$ws.Range("D18").Value = "'25.159.35"
$ws.Range("E18").Value = "'  -2.72%  "

# Row 19
$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "'  -2.96%  "

# Row 20
$ws.Range("D20").Value = "'0.000006516"
$ws.Range("E20").Value = "'  -3.61%  "

# Row 21
$ws.Range("D21").Value = "'4.415"
$ws.Range("E21").Value = "'  +3.36%  "

# Row 22
$ws.Range("D22").Value = "'1.881.59"
$ws.Range("E22").Value = "'  -4.28%  "

# Row 23
$ws.Range("D23").Value = "'8.436"
$ws.Range("E23").Value = "'  -2.15%  "

# Row 24
$ws.Range("D24").Value = "'5.259"
$ws.Range("E24").Value = "'  -0.29%  "

# Row 25
$ws.Range("D25").Value = "'132.79"
$ws.Range("E25").Value = "'  -4.32%  "

# Row 26
$ws.Range("D26").Value = "'14.69"
$ws.Range("E26").Value = "'  -3.13%  "

# Row 27
$ws.Range("D27").Value = "'1.384"
$ws.Range("E27").Value = "'  -8.70%  "

# Row 28
$ws.Range("D28").Value = "'1.700"
$ws.Range("E28").Value = "'  -3.62%  "

# Row 29
$ws.Range("D29").Value = "'102.29"
$ws.Range("E29").Value = "'  -3.10%  "

# Row 30
$ws.Range("D30").Value = "'3.805"
$ws.Range("E30").Value = "'  -3.36%  "

# Row 31
$ws.Range("D31").Value = "'0.07871"
$ws.Range("E31").Value = "'  -4.84%  "

# Row 32
$ws.Range("D32").Value = "'3.518"
$ws.Range("E32").Value = "'  -4.55%  "

# Row 33
$ws.Range("D33").Value = "'0.04586"
$ws.Range("E33").Value = "'  -0.87%  "

# Row 34
$ws.Range("D34").Value = "'2.613"
$ws.Range("E34").Value = "'  -1.13%  "

# Row 35
$ws.Range("D35").Value = "'0.9382"
$ws.Range("E35").Value = "'  -4.90%  "

# Row 36
$ws.Range("D36").Value = "'0.5827"
$ws.Range("E36").Value = "'  -5.70%  "

# Row 37
$ws.Range("D37").Value = "'2.611"
$ws.Range("E37").Value = "'  -2.62%  "

# Row 38
$ws.Range("D38").Value = "'0.8465"
$ws.Range("E38").Value = "'  +14.17%  "

# Row 39
$ws.Range("D39").Value = "'0.01538"
$ws.Range("E39").Value = "'  -3.95%  "

# Row 40
$ws.Range("D40").Value = "'1.002"
$ws.Range("E40").Value = "'  +0.23%  "

# Row 41
$ws.Range("D41").Value = "'1.834"
$ws.Range("E41").Value = "'  -5.33%  "

# Row 42
$ws.Range("D42").Value = "'98.21"
$ws.Range("E42").Value = "'  -0.05%  "

# Row 43
$ws.Range("D43").Value = "'0.3681"
$ws.Range("E43").Value = "'  -4.28%  "

# Row 44
$ws.Range("D44").Value = "'4.835"
$ws.Range("E44").Value = "'  -3.26%  "

# Row 45
$ws.Range("D45").Value = "'0.1132"
$ws.Range("E45").Value = "'  +0.22%  "

# Row 46
$ws.Range("B46").Value = "'Cronos"
$ws.Range("C46").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.05157"
$ws.Range("E46").Value = "'  -1.65%  "

# Row 47
$ws.Range("B47").Value = "'Aptos"
$ws.Range("C47").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'6.017"
$ws.Range("E47").Value = "'  -3.33%  "

# Row 48
$ws.Range("D48").Value = "'51.90"
$ws.Range("E48").Value = "'  -5.15%  "

# Row 49
$ws.Range("D49").Value = "'29.44"
$ws.Range("E49").Value = "'  -3.28%  "

# Row 50
$ws.Range("D50").Value = "'1.004"
$ws.Range("E50").Value = "'  +0.33%  "

# Row 51
$ws.Range("D51").Value = "'7.329"
$ws.Range("E51").Value = "'  -3.71%  "
